$d = $word.ActiveDocument
Write-Output ("StoryRanges.Count=" + $d.StoryRanges.Count)
$sr = $d.StoryRanges.Item(1)
Write-Output ("Story1 Text=[" + $sr.Text + "]")
$n = $sr.NextStoryRange
if ($n -ne $null) {
    Write-Output ("Story2 Text=[" + $n.Text + "]")
}
